$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "292.43"
$ws.Range("B2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-6.94%"
$ws.Range("B2").Copy()
$ws.Range("E2").PasteSpecial(-4122)

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.37"
$ws.Range("B3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.06%"
$ws.Range("B3").Copy()
$ws.Range("E3").PasteSpecial(-4122)

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.014"
$ws.Range("B4").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.90%"
$ws.Range("B4").Copy()
$ws.Range("E4").PasteSpecial(-4122)

# Row 5
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-3.44%"
$ws.Range("B5").Copy()
$ws.Range("E5").PasteSpecial(-4122)

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.529"
$ws.Range("B6").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-8.44%"
$ws.Range("B6").Copy()
$ws.Range("E6").PasteSpecial(-4122)

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9272"
$ws.Range("B7").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.13%"
$ws.Range("B7").Copy()
$ws.Range("E7").PasteSpecial(-4122)

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1193"
$ws.Range("B9").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.66%"
$ws.Range("B9").Copy()
$ws.Range("E9").PasteSpecial(-4122)

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1747"
$ws.Range("B10").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-3.86%"
$ws.Range("B10").Copy()
$ws.Range("E10").PasteSpecial(-4122)

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.04317"
$ws.Range("B11").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "5.09%"
$ws.Range("B11").Copy()
$ws.Range("E11").PasteSpecial(-4122)

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08639"
$ws.Range("B12").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.96%"
$ws.Range("B12").Copy()
$ws.Range("E12").PasteSpecial(-4122)

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1053"
$ws.Range("B13").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.05%"
$ws.Range("B13").Copy()
$ws.Range("E13").PasteSpecial(-4122)

# Row 14
$ws.Range("B14").Value = "TigerCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.006003"
$ws.Range("B14").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.01%"
$ws.Range("B14").Copy()
$ws.Range("E14").PasteSpecial(-4122)

# Row 15
$ws.Range("B15").Value = "LEO"
$ws.Range("C15").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.341"
$ws.Range("B15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.29%"
$ws.Range("B15").Copy()
$ws.Range("E15").PasteSpecial(-4122)

# Row 16
$ws.Range("B16").Value = "GateToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.299"
$ws.Range("B16").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.49%"
$ws.Range("B16").Copy()
$ws.Range("E16").PasteSpecial(-4122)

# Row 17
$ws.Range("B17").Value = "BitpandaEcosystemToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.3289"
$ws.Range("B17").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.02%"
$ws.Range("B17").Copy()
$ws.Range("E17").PasteSpecial(-4122)

# Row 18
$ws.Range("B18").Value = "MCDex"
$ws.Range("C18").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.975"
$ws.Range("B18").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "4.63%"
$ws.Range("B18").Copy()
$ws.Range("E18").PasteSpecial(-4122)

# Row 19
$ws.Range("B19").Value = "ProBitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.1390"
$ws.Range("B19").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.85%"
$ws.Range("B19").Copy()
$ws.Range("E19").PasteSpecial(-4122)

# Row 20
$ws.Range("B20").Value = "ZBToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.2795"
$ws.Range("B20").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.59%"
$ws.Range("B20").Copy()
$ws.Range("E20").PasteSpecial(-4122)

# Row 21
$ws.Range("B21").Value = "BitForexToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001280"
$ws.Range("B21").Copy()
$ws.Range("D21").PasteSpecial(-4122)
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.32%"
$ws.Range("B21").Copy()
$ws.Range("E21").PasteSpecial(-4122)

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.03943"
$ws.Range("B22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.11%"
$ws.Range("B22").Copy()
$ws.Range("E22").PasteSpecial(-4122)

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.93%"
$ws.Range("B23").Copy()
$ws.Range("E23").PasteSpecial(-4122)

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.003778"
$ws.Range("B24").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-4.98%"
$ws.Range("B24").Copy()
$ws.Range("E24").PasteSpecial(-4122)

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.92%"
$ws.Range("B25").Copy()
$ws.Range("E25").PasteSpecial(-4122)

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003725"
$ws.Range("B26").Copy()
$ws.Range("D26").PasteSpecial(-4122)

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02284"
$ws.Range("B38").Copy()
$ws.Range("D38").PasteSpecial(-4122)
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-5.65%"
$ws.Range("B38").Copy()
$ws.Range("E38").PasteSpecial(-4122)

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.04976"
$ws.Range("B39").Copy()
$ws.Range("D39").PasteSpecial(-4122)
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-3.53%"
$ws.Range("B39").Copy()
$ws.Range("E39").PasteSpecial(-4122)

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.005633"
$ws.Range("B40").Copy()
$ws.Range("D40").PasteSpecial(-4122)
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "70.72%"
$ws.Range("B40").Copy()
$ws.Range("E40").PasteSpecial(-4122)

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007700"
$ws.Range("B41").Copy()
$ws.Range("D41").PasteSpecial(-4122)
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.05%"
$ws.Range("B41").Copy()
$ws.Range("E41").PasteSpecial(-4122)

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1285"
$ws.Range("B42").Copy()
$ws.Range("D42").PasteSpecial(-4122)

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.61%"
$ws.Range("B43").Copy()
$ws.Range("E43").PasteSpecial(-4122)

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007299"
$ws.Range("B44").Copy()
$ws.Range("D44").PasteSpecial(-4122)
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-11.46%"
$ws.Range("B44").Copy()
$ws.Range("E44").PasteSpecial(-4122)

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.2923"
$ws.Range("B45").Copy()
$ws.Range("D45").PasteSpecial(-4122)
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-13.90%"
$ws.Range("B45").Copy()
$ws.Range("E45").PasteSpecial(-4122)

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006324"
$ws.Range("B46").Copy()
$ws.Range("D46").PasteSpecial(-4122)
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.06%"
$ws.Range("B46").Copy()
$ws.Range("E46").PasteSpecial(-4122)

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.11%"
$ws.Range("B47").Copy()
$ws.Range("E47").PasteSpecial(-4122)

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.02132"
$ws.Range("B48").Copy()
$ws.Range("D48").PasteSpecial(-4122)
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-72.22%"
$ws.Range("B48").Copy()
$ws.Range("E48").PasteSpecial(-4122)

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.11%"
$ws.Range("B49").Copy()
$ws.Range("E49").PasteSpecial(-4122)

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.11%"
$ws.Range("B50").Copy()
$ws.Range("E50").PasteSpecial(-4122)
